$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right after the current row 110 (old "Vega Central Mapocho
# de Santiago" / Camote record dated 44263). This pushes the existing row 111
# (dated 44845) down to row 113, and leaves two blank rows (111, 112) to fill in.
$ws.Rows.Item(111).Insert()
$ws.Rows.Item(111).Insert()

# Row 110 is updated in place with a new weekly price report (2022-10-24).
$ws.Cells.Item(110, 1).Value = 9
$ws.Cells.Item(110, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(110, 3).Value = "Metropolitana"
$ws.Cells.Item(110, 4).Value = 44858
$ws.Cells.Item(110, 5).Value = 13
$ws.Cells.Item(110, 6).Value = 100114002
$ws.Cells.Item(110, 7).Value = "Camote"
$ws.Cells.Item(110, 8).Value = "Sin especificar"
$ws.Cells.Item(110, 9).Value = "Primera"
$ws.Cells.Item(110, 10).Value = 750
$ws.Cells.Item(110, 11).Value = 17000
$ws.Cells.Item(110, 12).Value = 18000
$ws.Cells.Item(110, 13).Value = 17533
$ws.Cells.Item(110, 14).Value = "`$/caja 18 kilos"
$ws.Cells.Item(110, 15).Value = "Perú"
$ws.Cells.Item(110, 16).Value = 974
$ws.Cells.Item(110, 17).Value = 18
$ws.Cells.Item(110, 18).Value = "Hortaliza"

# Row 111 is a brand new record, also dated 2022-10-24.
$ws.Cells.Item(111, 1).Value = 9
$ws.Cells.Item(111, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(111, 3).Value = "Metropolitana"
$ws.Cells.Item(111, 4).Value = 44858
$ws.Cells.Item(111, 5).Value = 13
$ws.Cells.Item(111, 6).Value = 100114002
$ws.Cells.Item(111, 7).Value = "Camote"
$ws.Cells.Item(111, 8).Value = "Sin especificar"
$ws.Cells.Item(111, 9).Value = "Primera"
$ws.Cells.Item(111, 10).Value = 1500
$ws.Cells.Item(111, 11).Value = 13000
$ws.Cells.Item(111, 12).Value = 14000
$ws.Cells.Item(111, 13).Value = 13533
$ws.Cells.Item(111, 14).Value = "`$/malla 18 kilos"
$ws.Cells.Item(111, 15).Value = "Perú"
$ws.Cells.Item(111, 16).Value = 752
$ws.Cells.Item(111, 17).Value = 18
$ws.Cells.Item(111, 18).Value = "Hortaliza"

# Row 112 restores the original row-110 record (the one that used to sit at 110
# before today's update), now pushed down below the new entries.
$ws.Cells.Item(112, 1).Value = 9
$ws.Cells.Item(112, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(112, 3).Value = "Metropolitana"
$ws.Cells.Item(112, 4).Value = 44263
$ws.Cells.Item(112, 5).Value = 13
$ws.Cells.Item(112, 6).Value = 100114002
$ws.Cells.Item(112, 7).Value = "Camote"
$ws.Cells.Item(112, 8).Value = "Sin especificar"
$ws.Cells.Item(112, 9).Value = "Primera"
$ws.Cells.Item(112, 10).Value = 1600
$ws.Cells.Item(112, 11).Value = 9000
$ws.Cells.Item(112, 12).Value = 9000
$ws.Cells.Item(112, 13).Value = 9000
$ws.Cells.Item(112, 14).Value = "`$/malla 18 kilos"
$ws.Cells.Item(112, 15).Value = "Perú"
$ws.Cells.Item(112, 16).Value = 500
$ws.Cells.Item(112, 17).Value = 18
$ws.Cells.Item(112, 18).Value = "Hortaliza"

# Row 113 already holds the shifted-down former row 111 values (date 44845,
# etc.) thanks to the row insert above, so nothing further to do there.

# Column D carries a date/time number format (style index 2 in the original
# file) -- make sure it stays applied on the two brand-new rows we just wrote,
# by copying the format from the still-correctly-styled row 113 (originally
# row 111) down onto rows 110-112.
$ws.Range("D113").Copy()
$ws.Range("D110:D112").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Keep the sheet's used-range / dimension declaration (A1:R113) accurate.
$ws.Range("A1").Select()
